{"js": "// Append two new paragraphs after the existing one:\n//   1) an empty paragraph\n//   2) \"But recently we changed something in footer and header! Stay update!\"\n// The \"_GoBack\" bookmark originally sitting at the end of the first\n// paragraph is moved so it still sits right after \"...header! Sta\",\n// splitting \"Stay update!\" into \"Sta\" + \"y update!\" (matching the target\n// diff's two runs around the bookmark).\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = body.paragraphs.items[0];\nconst tailOfFirstParagraph = firstParagraph.getRange(\"End\");\n\n// Insert a genuinely empty paragraph (no run) right after paragraph 1,\n// matching what Word produces for a blank line whose mark only carries\n// run-formatting (<w:pPr><w:rPr>...). Plain insertParagraph()/insertText()\n// with \"\\n\" leaves a stray empty run behind, so we splice in the exact\n// OOXML for the blank paragraph instead.\nconst blankParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\ntailOfFirstParagraph.insertOoxml(blankParagraphOoxml, Word.InsertLocation.end);\nawait context.sync();\n\n// Now append the new sentence as its own paragraph after that blank line.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst blankParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\nconst tailOfBlankParagraph = blankParagraph.getRange(\"End\");\ntailOfBlankParagraph.insertText(\n  \"\\nBut recently we changed something in footer and header! Stay update!\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\n// The \"_GoBack\" bookmark stays attached to the end of paragraph 1 after the\n// insertions above; remove it from there and re-create it in its new\n// location (right before \"y update!\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst found = body.search(\n  \"But recently we changed something in footer and header! Sta\",\n  { matchCase: true, matchWholeWord: false }\n);\nfound.load(\"items\");\nawait context.sync();\n\nconst targetRange = found.items[0].getRange(\"End\");\ntargetRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Append two new paragraphs after the existing one:\n#   1) an empty paragraph\n#   2) \"But recently we changed something in footer and header! Stay update!\"\n# The \"_GoBack\" bookmark originally sitting at the end of the first\n# paragraph is moved so it still sits right after \"...header! Sta\",\n# splitting \"Stay update!\" into \"Sta\" + \"y update!\" (matching the target\n# diff's two runs around the bookmark).\n\n$d = $word.ActiveDocument\n\n# Insert a genuinely empty paragraph (no run) at the very end of the\n# document, matching what Word produces for a blank line whose mark only\n# carries run-formatting (<w:pPr><w:rPr>...). A plain\n# InsertParagraphAfter()/InsertAfter(\"`r\") leaves a stray empty run behind,\n# so splice in the exact OOXML for the blank paragraph instead.\n$endOfDoc = $d.Content\n$endOfDoc.Collapse(0)   # wdCollapseEnd\n$blankParagraphXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$endOfDoc.InsertXML($blankParagraphXml)\n$d.Save()\n\n# Now append the new sentence as its own paragraph after that blank line.\n$tail = $d.Content\n$tail.Collapse(0)   # wdCollapseEnd\n$tail.InsertAfter(\"`rBut recently we changed something in footer and header! Stay update!\")\n$d.Save()\n\n# Move the \"_GoBack\" bookmark from the end of paragraph 1 to its new home.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n$found = $d.Content\n$found.Find.Execute(\"But recently we changed something in footer and header! Sta\")\n$found.Collapse(0)   # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $found)\n$d.Save()\n"}
